$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-valued columns (D = Price, B/C/E) keep their exact string
# representation instead of being auto-coerced to numbers/dates by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.338.90"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.59%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.375.09"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -4.00%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.37%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "501.54"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.37"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.63%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.544"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.381.86"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.84%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0983"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.58%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.55%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.325"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.65"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.79%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.796.78"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.78%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "56.292.36"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.76%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.58"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.61%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.18%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.380.20"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.94%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.06"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.01"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.45%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "307.29"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.24"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.64%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.00"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.36%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.368"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.69%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.94%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.28"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -5.35%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "172.49"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.00%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0714"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.94%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.11%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.24%  "

$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.76"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -8.16%  "

$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.21%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.08"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.92%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.60"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.90%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -7.04%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.78"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.52%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.02"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.58%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.794"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.41%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.41"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.83%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "131.16"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.71%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.34"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.86%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.77"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.96%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.567"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.90%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0904"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.01%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "242.39"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -7.76%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0483"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.73%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0209"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.43%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.00"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.44%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.55"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.30%  "
